# Generate Report for Archive
# Reorders rows 7-9 in Overview/zh-cn/de-de sheets: e8088aa3 (now 'In Translation')
# moves to row 7, pushing 96e8afa5 and cf667357 down by one row each.
$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item(1)
$ws.Range('A7').Value = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md'
$ws.Range('B7').Value = 'In Translation'
$ws.Range('C7').Value = 'In Translation'
$ws.Range('D7').Value = '2016-03-22 00:38:57'
$ws.Range('A8').Value = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md'
$ws.Range('B8').Value = 'Ready for handoff'
$ws.Range('C8').Value = 'Ready for handoff'
$ws.Range('D8').Value = '2016-03-22 00:33:37'
$ws.Range('A9').Value = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.md'
$ws.Range('B9').Value = 'Ready for handoff'
$ws.Range('C9').Value = 'Ready for handoff'
$ws.Range('D9').Value = '2016-03-22 00:36:31'

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$7') { $hl.TextToDisplay = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md' }
    if ($addr -eq '$A$8') { $hl.TextToDisplay = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md' }
    if ($addr -eq '$A$9') { $hl.TextToDisplay = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.md' }
}

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item(2)
$ws.Range('A7').Value = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md'
$ws.Range('B7').Value = '.md'
$ws.Range('C7').Value = 'In Translation'
$ws.Range('D7').Value = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.zh-cn.xlf'
$ws.Range('E7').Value = '2016-03-22 00:38:54'
$ws.Range('H7').Value = '0001-01-01 00:00:00'
$ws.Range('J7').Value = 'Include'
$ws.Range('A8').Value = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md'
$ws.Range('B8').Value = '.md'
$ws.Range('C8').Value = 'Ready for handoff'
$ws.Range('D8').Value = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.zh-cn.xlf'
$ws.Range('E8').Value = '2016-03-22 00:33:33'
$ws.Range('H8').Value = '0001-01-01 00:00:00'
$ws.Range('J8').Value = 'Include'
$ws.Range('A9').Value = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.md'
$ws.Range('B9').Value = '.md'
$ws.Range('C9').Value = 'Ready for handoff'
$ws.Range('D9').Value = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.bb122b9ccdade679a9783ff69492a289cd8dd1fb.zh-cn.xlf'
$ws.Range('E9').Value = '2016-03-22 00:36:27'
$ws.Range('H9').Value = '0001-01-01 00:00:00'
$ws.Range('J9').Value = 'Include'

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$7') { $hl.TextToDisplay = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md' }
    if ($addr -eq '$D$7') { $hl.TextToDisplay = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.zh-cn.xlf' }
    if ($addr -eq '$A$8') { $hl.TextToDisplay = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md' }
    if ($addr -eq '$D$8') { $hl.TextToDisplay = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.zh-cn.xlf' }
    if ($addr -eq '$A$9') { $hl.TextToDisplay = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.md' }
    if ($addr -eq '$D$9') { $hl.TextToDisplay = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.bb122b9ccdade679a9783ff69492a289cd8dd1fb.zh-cn.xlf' }
}

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item(3)
$ws.Range('A7').Value = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md'
$ws.Range('B7').Value = '.md'
$ws.Range('C7').Value = 'In Translation'
$ws.Range('D7').Value = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.de-de.xlf'
$ws.Range('E7').Value = '2016-03-22 00:38:57'
$ws.Range('H7').Value = '0001-01-01 00:00:00'
$ws.Range('J7').Value = 'Include'
$ws.Range('A8').Value = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md'
$ws.Range('B8').Value = '.md'
$ws.Range('C8').Value = 'Ready for handoff'
$ws.Range('D8').Value = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.de-de.xlf'
$ws.Range('E8').Value = '2016-03-22 00:33:37'
$ws.Range('H8').Value = '0001-01-01 00:00:00'
$ws.Range('J8').Value = 'Include'
$ws.Range('A9').Value = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.md'
$ws.Range('B9').Value = '.md'
$ws.Range('C9').Value = 'Ready for handoff'
$ws.Range('D9').Value = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.bb122b9ccdade679a9783ff69492a289cd8dd1fb.de-de.xlf'
$ws.Range('E9').Value = '2016-03-22 00:36:31'
$ws.Range('H9').Value = '0001-01-01 00:00:00'
$ws.Range('J9').Value = 'Include'

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$7') { $hl.TextToDisplay = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.md' }
    if ($addr -eq '$D$7') { $hl.TextToDisplay = 'e8088aa3-5fd2-41a5-8060-ea3b75c18b96.e1568a30eeff22474690f39448d8625f901cd9e9.de-de.xlf' }
    if ($addr -eq '$A$8') { $hl.TextToDisplay = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.md' }
    if ($addr -eq '$D$8') { $hl.TextToDisplay = '96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8.0012e40d796e5c6f54b3c87d5af7bf616b8ae37b.de-de.xlf' }
    if ($addr -eq '$A$9') { $hl.TextToDisplay = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.md' }
    if ($addr -eq '$D$9') { $hl.TextToDisplay = 'cf667357-71fc-4ed0-8b1d-8d1dc74b41dd.bb122b9ccdade679a9783ff69492a289cd8dd1fb.de-de.xlf' }
}
